# Add a new "Greece" tab for the Greece Market test data, modeled on the
# existing "Croatia" sheet (same layout/styles, new market name + ticket ref).

$wb = $excel.ActiveWorkbook
$croatia = $wb.Worksheets.Item("Croatia")

# Record a "select-all" state on the Croatia sheet (it stops being the active
# tab once the new sheet takes over).
$croatia.Range("A1:XFD1048576").Select() | Out-Null

# Duplicate Croatia right after itself to become the new Greece sheet.
$croatia.Copy([System.Reflection.Missing]::Value, $croatia)
$greece = $wb.Worksheets.Item($croatia.Index + 1)
$greece.Name = "Greece"

# Fill in the Greece-specific values (ticket ref first, then market name, to
# match shared-string insertion order).
$greece.Range("B4").Value = "NGC-4119/T3167"
$greece.Range("B2").Value = "Greece Market"

# Leave the new sheet selected/active, matching the recorded cursor position.
$greece.Range("H23").Select() | Out-Null
